# Wrapping up DEB manuscript and submitting it
#
# This script reproduces the author's final proofreading pass: several
# checklist bullets get a grammar-checker mark (w:proofErr gramStart/gramEnd)
# wrapped around one word (this is what Word's "Check Grammar" feature
# stamps into the OOXML when it flags a phrase, splitting the run that
# contained it into three runs), and one bullet gets a genuine wording fix
# ("Explain in significance" -> "Explain significance").

$d = $word.ActiveDocument

function Get-ParaByText {
    param([object]$doc, [string]$searchText)
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.TrimEnd([char]13) -eq $searchText) {
            return $p
        }
    }
    return $null
}

function Set-ParaRunsXml {
    param([object]$para, [string]$paraAttrs, [string]$pPrXml, [string]$runsXml)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body><w:p ' + $paraAttrs + '>' + $pPrXml + $runsXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml) | Out-Null
}

function Split-ParaWithGram {
    param([object]$doc, [string]$oldText, [string]$paraAttrs, [string]$pPrXml, [string[]]$segments)
    # $segments alternates: plain-text, gram-word, plain-text, gram-word, ...
    # (always an odd count, starting and ending with a plain-text segment,
    # each gram-word wrapped in proofErr gramStart/gramEnd).
    $para = Get-ParaByText $doc $oldText
    if ($para -eq $null) {
        throw "Paragraph not found: $oldText"
    }

    $runs = ""
    for ($i = 0; $i -lt $segments.Length; $i++) {
        $seg = $segments[$i]
        $escaped = $seg.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        if ($i % 2 -eq 0) {
            $runs += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
        } else {
            $runs += '<w:proofErr w:type="gramStart"/><w:r><w:t>' + $escaped + '</w:t></w:r><w:proofErr w:type="gramEnd"/>'
        }
    }

    Set-ParaRunsXml $para $paraAttrs $pPrXml $runs
}

$pPrLvl1 = '<w:pPr><w:pStyle w:val="TS"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$pPrLvl2 = '<w:pPr><w:pStyle w:val="TS"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# 1) "Indicate if color should be used for any figures in print"
Split-ParaWithGram $d `
    "Indicate if color should be used for any figures in print" `
    'w14:paraId="137133D5" w14:textId="331C56C2" w:rsidR="008057A5" w:rsidRDefault="008057A5" w:rsidP="008057A5"' `
    $pPrLvl1 `
    @("Indicate if color should be used for any figures in ", "print")

# 2) "Check reference list/in text citations, figure and table numbers, spelling and grammar. "
Split-ParaWithGram $d `
    "Check reference list/in text citations, figure and table numbers, spelling and grammar. " `
    'w14:paraId="45250178" w14:textId="60B69525" w:rsidR="008057A5" w:rsidRDefault="008057A5" w:rsidP="008057A5"' `
    $pPrLvl1 `
    @("Check reference list/in text citations, figure and table numbers, ", "spelling", " and grammar. ")

# 3) "Competing interests statement. "
Split-ParaWithGram $d `
    "Competing interests statement. " `
    'w14:paraId="1350FC1F" w14:textId="78FA78C5" w:rsidR="008057A5" w:rsidRDefault="008057A5" w:rsidP="008057A5"' `
    $pPrLvl1 `
    @("Competing ", "interests", " statement. ")

# 4) "Has the model been presented in such details that the reader is able to develop the model?"
Split-ParaWithGram $d `
    "Has the model been presented in such details that the reader is able to develop the model?" `
    'w14:paraId="7D6173EA" w14:textId="3755204D" w:rsidR="008057A5" w:rsidRDefault="008057A5" w:rsidP="008057A5"' `
    $pPrLvl2 `
    @("Has the model been presented in such ", "details", " that the reader is able to develop the model?")

# 5) "Can put everything in one word/pdf file, or include figures as source files. "
Split-ParaWithGram $d `
    "Can put everything in one word/pdf file, or include figures as source files. " `
    'w14:paraId="541CBA58" w14:textId="3B0BD9FF" w:rsidR="00192B15" w:rsidRDefault="00192B15" w:rsidP="00192B15"' `
    $pPrLvl1 `
    @("Can put everything in one word/pdf ", "file, or", " include figures as source files. ")

# 6) "Essential elements: ... conclusions, artwork and tables with captions. "
Split-ParaWithGram $d `
    "Essential elements: abstract, keywords, introduction, materials and methods, results, conclusions, artwork and tables with captions. " `
    'w14:paraId="0D9C2E26" w14:textId="0247D33C" w:rsidR="00192B15" w:rsidRDefault="00192B15" w:rsidP="00192B15"' `
    $pPrLvl1 `
    @("Essential elements: abstract, keywords, introduction, materials and methods, results, conclusions, ", "artwork", " and tables with captions. ")

# 7) Discussion bullet: genuine wording correction, no proofErr marks
$d.Content.Find.Execute(
    "Explain in significance in the context of the current literature.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Explain significance in the context of the current literature.", 2) | Out-Null

# 8) "Conclusions: need to add this. Or use existing last 2 paragraphs. ..."
Split-ParaWithGram $d `
    "Conclusions: need to add this. Or use existing last 2 paragraphs. Did dissertation chapter have conclusions section?" `
    'w14:paraId="4A84B728" w14:textId="762ACF31" w:rsidR="00192B15" w:rsidRDefault="00192B15" w:rsidP="00192B15"' `
    $pPrLvl1 `
    @("Conclusions: need to add this. Or use ", "existing", " last 2 paragraphs. Did dissertation chapter have conclusions section?")
